$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sample")

# Update data rows (row 2 - John/J/Doe - stays unchanged)
$ws.Range("A3").Value = "Katie"
$ws.Range("C3").Value = "Ball"

$ws.Range("A4").Value = "Donald"
$ws.Range("C4").Value = "Trump"

$ws.Range("A5").Value = "Mohammed"
$ws.Range("C5").Value = "Salah"

# Update the active selection to C6
$ws.Range("C6").Select()
